$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.977.14"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3
$ws.Range("D3").Value = "2.263.70"
$ws.Range("E3").Value = "  +1.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "318.65"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6
$ws.Range("D6").Value = "101.87"
$ws.Range("E6").Value = "  +1.21%  "

# Row 7
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -1.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  -1.44%  "

# Row 10
$ws.Range("D10").Value = "37.26"

# Row 11
$ws.Range("D11").Value = "0.0833"
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  -0.54%  "

# Row 13
$ws.Range("E13").Value = "  -2.27%  "

# Row 14
$ws.Range("D14").Value = "2.604.04"
$ws.Range("E14").Value = "  +0.83%  "

# Row 15
$ws.Range("D15").Value = "0.864"
$ws.Range("E15").Value = "  -0.62%  "

# Row 16
$ws.Range("D16").Value = "14.45"
$ws.Range("E16").Value = "  -0.21%  "

# Row 17
$ws.Range("D17").Value = "2.259.81"
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$ws.Range("D18").Value = "43.902.87"
$ws.Range("E18").Value = "  +1.66%  "

# Row 19
$ws.Range("D19").Value = "13.36"
$ws.Range("E19").Value = "  -8.81%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +1.43%  "

# Row 21
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("D22").Value = "65.77"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
$ws.Range("E23").Value = "  -1.65%  "

# Row 24
$ws.Range("D24").Value = "235.95"
$ws.Range("E24").Value = "  -1.13%  "

# Row 25
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  -2.94%  "

# Row 26
$ws.Range("E26").Value = "  -0.16%  "

# Row 27
$ws.Range("D27").Value = "10.16"
$ws.Range("E27").Value = "  +1.12%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.17"
$ws.Range("E28").Value = "  -3.68%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "37.43"
$ws.Range("E29").Value = "  +2.90%  "

# Row 30
$ws.Range("D30").Value = "6.24"
$ws.Range("E30").Value = "  -2.07%  "

# Row 31
$ws.Range("D31").Value = "20.27"
$ws.Range("E31").Value = "  -0.92%  "

# Row 32
$ws.Range("D32").Value = "157.79"
$ws.Range("E32").Value = "  +4.86%  "

# Row 33
$ws.Range("E33").Value = "  -3.14%  "

# Row 34
$ws.Range("E34").Value = "  +0.17%  "

# Row 35
$ws.Range("D35").Value = "0.115"
$ws.Range("E35").Value = "  +10.34%  "

# Row 36
$ws.Range("E36").Value = "  -3.34%  "

# Row 37
$ws.Range("E37").Value = "  +0.33%  "

# Row 38
$ws.Range("E38").Value = "  -2.45%  "

# Row 39
$ws.Range("D39").Value = "16.07"
$ws.Range("E39").Value = "  +17.50%  "

# Row 40
$ws.Range("D40").Value = "3.72"
$ws.Range("E40").Value = "  +1.30%  "

# Row 41
$ws.Range("D41").Value = "'4.20"
$ws.Range("E41").Value = "  -5.80%  "

# Row 42
$ws.Range("E42").Value = "  -2.57%  "

# Row 43
$ws.Range("E43").Value = "  +0.11%  "

# Row 44
$ws.Range("D44").Value = "1.797.45"

# Row 45
$ws.Range("E45").Value = "  -3.03%  "

# Row 46
$ws.Range("D46").Value = "75.75"
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").Value = "82.53"
$ws.Range("E47").Value = "  -4.95%  "

# Row 48
$ws.Range("D48").Value = "5.22"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "104.36"
$ws.Range("E49").Value = "  +1.03%  "

# Row 50
$ws.Range("D50").Value = "58.62"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.66"
$ws.Range("E51").Value = "  +4.49%  "
